# Auto-generated Excel COM-interop edit script
# Applies updated Leve profit figures (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# matching the "chore: update Sheets via scheduled runner" data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 695.4627
$ws.Range("I15").Value = 695.4627
$ws.Range("K15").Value = 2086.3881
$ws.Range("M15").Value = -1917.3881
$ws.Range("H33").Value = 341.4
$ws.Range("I33").Value = 186.33333
$ws.Range("K33").Value = 186.33333
$ws.Range("M33").Value = 42.66667000000001
$ws.Range("H88").Value = 8329.666999999999
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 7494.5
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 7494.5
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -8306.5
$ws.Range("H91").Value = 8329.666999999999
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 7494.5
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 7494.5
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -10302.5
$ws.Range("H98").Value = 1559.0834
$ws.Range("I98").Value = 1473.5454
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 1473.5454
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = 24.45460000000003
$ws.Range("N98").Value = -5496
$ws.Range("H108").Value = 47358
$ws.Range("J108").Value = 47644
$ws.Range("L108").Value = 47644
$ws.Range("N108").Value = -55324
$ws.Range("H122").Value = 1559.0834
$ws.Range("I122").Value = 1473.5454
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4420.6362
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1970.6362
$ws.Range("N122").Value = -12400
$ws.Range("H135").Value = 2247.75
$ws.Range("I135").Value = 900
$ws.Range("K135").Value = 8100
$ws.Range("M135").Value = -5565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7438.1304
$ws.Range("I32").Value = 7438.1304
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7438.1304
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7151.1304
$ws.Range("N32").ClearContents()
$ws.Range("H88").Value = 15930.286
$ws.Range("I88").Value = 34204
$ws.Range("J88").Value = 2225
$ws.Range("K88").Value = 34204
$ws.Range("L88").Value = 2225
$ws.Range("M88").Value = -33798
$ws.Range("N88").Value = -3037
$ws.Range("H91").Value = 15930.286
$ws.Range("I91").Value = 34204
$ws.Range("J91").Value = 2225
$ws.Range("K91").Value = 34204
$ws.Range("L91").Value = 2225
$ws.Range("M91").Value = -32800
$ws.Range("N91").Value = -5033
$ws.Range("H97").Value = 1047.7632
$ws.Range("I97").Value = 1059.2174
$ws.Range("K97").Value = 1059.2174
$ws.Range("M97").Value = -563.2174
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 203816.3
$ws.Range("I102").Value = 336162.84
$ws.Range("K102").Value = 336162.84
$ws.Range("M102").Value = -334540.84
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840
$ws.Range("H132").Value = 29461542
$ws.Range("I132").Value = 9835.714
$ws.Range("K132").Value = 29507.142
$ws.Range("M132").Value = -26977.142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10872.454
$ws.Range("I105").Value = 13849.625
$ws.Range("J105").Value = 2933.3333
$ws.Range("K105").Value = 13849.625
$ws.Range("L105").Value = 2933.3333
$ws.Range("M105").Value = -12102.625
$ws.Range("N105").Value = -6427.3333
$ws.Range("H107").Value = 3308.4666
$ws.Range("I107").Value = 2779.5386
$ws.Range("K107").Value = 2779.5386
$ws.Range("M107").Value = -859.5385999999999
$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82774

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3189.3684
$ws.Range("I31").Value = 2135.1707
$ws.Range("J31").Value = 5890.75
$ws.Range("K31").Value = 2135.1707
$ws.Range("L31").Value = 5890.75
$ws.Range("M31").Value = -1840.1707
$ws.Range("N31").Value = -6480.75
$ws.Range("H34").Value = 3189.3684
$ws.Range("I34").Value = 2135.1707
$ws.Range("J34").Value = 5890.75
$ws.Range("K34").Value = 2135.1707
$ws.Range("L34").Value = 5890.75
$ws.Range("M34").Value = -1933.1707
$ws.Range("N34").Value = -6294.75
$ws.Range("H58").Value = 3269.5454
$ws.Range("I58").Value = 2933.125
$ws.Range("K58").Value = 2933.125
$ws.Range("M58").Value = -2730.125
$ws.Range("H122").Value = 1528.25
$ws.Range("I122").Value = 1350.7894
$ws.Range("K122").Value = 4052.3682
$ws.Range("M122").Value = -1602.3682
$ws.Range("H132").Value = 47434.18
$ws.Range("I132").Value = 55230.105
$ws.Range("J132").Value = 5113.4287
$ws.Range("K132").Value = 165690.315
$ws.Range("L132").Value = 15340.2861
$ws.Range("M132").Value = -163160.315
$ws.Range("N132").Value = -20400.2861
$ws.Range("H134").Value = 1291.8276
$ws.Range("I134").Value = 1294.7307
$ws.Range("K134").Value = 3884.1921
$ws.Range("M134").Value = -1349.1921
$ws.Range("H135").Value = 55000
$ws.Range("J135").Value = 55000
$ws.Range("L135").Value = 55000
$ws.Range("N135").Value = -65140
$ws.Range("H136").Value = 3269.5454
$ws.Range("I136").Value = 2933.125
$ws.Range("K136").Value = 8799.375
$ws.Range("M136").Value = -6249.375
$ws.Range("H141").Value = 207832.58
$ws.Range("J141").Value = 207832.58
$ws.Range("L141").Value = 207832.58
$ws.Range("N141").Value = -218192.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 47600812
$ws.Range("I4").Value = 60142564
$ws.Range("K4").Value = 180427692
$ws.Range("M4").Value = -180427580
$ws.Range("H51").Value = 832.4286
$ws.Range("I51").Value = 365.4
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 1096.2
$ws.Range("L51").Value = 6000
$ws.Range("M51").Value = -636.1999999999998
$ws.Range("N51").Value = -6920
$ws.Range("H64").Value = 23299.8
$ws.Range("J64").Value = 28249.75
$ws.Range("L64").Value = 84749.25
$ws.Range("N64").Value = -85289.25
$ws.Range("H67").Value = 23299.8
$ws.Range("J67").Value = 28249.75
$ws.Range("L67").Value = 84749.25
$ws.Range("N67").Value = -86621.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2299.9092
$ws.Range("I132").Value = 2000.1666
$ws.Range("K132").Value = 6000.4998
$ws.Range("M132").Value = -3470.4998
$ws.Range("H139").Value = 124965.336
$ws.Range("J139").Value = 124965.336
$ws.Range("L139").Value = 124965.336
$ws.Range("N139").Value = -135245.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2121.093
$ws.Range("I61").Value = 1406.8077
$ws.Range("K61").Value = 1406.8077
$ws.Range("M61").Value = -1204.8077
$ws.Range("H68").Value = 6221.5
$ws.Range("I68").Value = 2776.3333
$ws.Range("K68").Value = 2776.3333
$ws.Range("M68").Value = -2027.3333
$ws.Range("H71").Value = 6221.5
$ws.Range("I71").Value = 2776.3333
$ws.Range("K71").Value = 13881.6665
$ws.Range("M71").Value = -10137.6665
$ws.Range("H113").Value = 2121.093
$ws.Range("I113").Value = 1406.8077
$ws.Range("K113").Value = 1406.8077
$ws.Range("M113").Value = 763.1922999999999
$ws.Range("H132").Value = 4709.4116
$ws.Range("I132").Value = 4070.6667
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 12212.0001
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -9682.000100000001
$ws.Range("N132").Value = -33560
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("H136").Value = 1711.85
$ws.Range("I136").Value = 1453
$ws.Range("K136").Value = 4359
$ws.Range("M136").Value = -1809
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3211.4666
$ws.Range("I132").Value = 3369.7856
$ws.Range("K132").Value = 10109.3568
$ws.Range("M132").Value = -7579.356800000001
$ws.Range("H136").Value = 1897.1666
$ws.Range("I136").Value = 1096.35
$ws.Range("K136").Value = 3289.05
$ws.Range("M136").Value = -739.0499999999997
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

